$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the two salary columns (O, P) and shift the remaining headers,
#     adding a new "Tiene cuenta" column in R -----------------------------
#   Old: N=Oficina | O=Salario Bruto | P=Salario Neto | Q=Fecha de ingreso empresa | R=Numero ID
#   New: N=Oficina | O=Salario Promedio | P=Fecha de ingreso empresa | Q=Numero ID | R=Tiene cuenta

$ws.Range("O1").Value = "Salario Promedio"
$ws.Range("P1").Value = "Fecha de ingreso empresa"
$ws.Range("Q1").Value = "Número ID"
$ws.Range("R1").Value = "Tiene cuenta"

# P1 becomes a "text style" header, matching the other text columns (e.g. A1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null

# Q1 keeps the "number style" header used by the other numeric-ish columns (e.g. N1)
$ws.Range("N1").Copy() | Out-Null
$ws.Range("Q1").PasteSpecial(-4122) | Out-Null

# R1 (new column) starts from the same "number style" header as N1 for the
# font/fill/number-format, then gets a left/right-only border (no top/bottom)
# to visually close the table, built up from nothing to avoid disturbing the
# colors of the existing shared full-border style.
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$ws.Range("R1").Borders.LineStyle = -4142
$ws.Range("R1").Borders.Item(7).LineStyle = 1
$ws.Range("R1").Borders.Item(10).LineStyle = 1

# --- Column widths: best-fit-ish widths for the shifted/new header columns --
$ws.Columns.Item(15).ColumnWidth = 15.4518
$ws.Columns.Item(16).ColumnWidth = 23.4518
$ws.Columns.Item(17).ColumnWidth = 12.3073
$ws.Columns.Item(18).ColumnWidth = 11.5924

# --- Selection / view state ------------------------------------------------
$ws.Range("Q8").Select() | Out-Null
